# Contributor 1 - ExampleData2
# Add 2 new columns (Age and BMI Status) to the "Data" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Header row
$ws.Range("D1").Value = "Age"
$ws.Range("E1").Value = "BMI Status"

# Age values (column D)
$ages = @(25, 34, 40, 31, 33, 22, 24, 27, 28, 34, 23, 19, 17, 55)
for ($i = 0; $i -lt $ages.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $ages[$i]
}

# BMI Status values (column E)
$statuses = @("Healthy", "Healthy", "NA", "Healthy", "Healthy", "NA", "Unhealthy", "Unhealthy", "Healthy", "NA", "NA", "Healthy", "Healthy", "Healthy")
for ($i = 0; $i -lt $statuses.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $statuses[$i]
}

# Update selection to reflect final cursor position (E12)
$ws.Range("E12").Select()
